$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update organization/contact info in column B (rows 6-10)
$ws.Range("B6").Value = "National Statistical Committee of the Kyrgyz Republic (Department of Household Statistics)"
$ws.Range("B7").Value = "Kalymbetova Yryskan"
$ws.Range("B8").Value = "yryskan.kalymbetova@gmail.com "
$ws.Range("B9").Value = "(0312) 32 46 55"
$ws.Range("B10").Value = "www.stat.gov.kg"

# Move the active selection to B4
$null = $ws.Range("B4").Select()
